$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column N/O (years 2021, 2022) added to the table ---

# Row 3: thin empty separator cells, same style as K3:M3 (s=29)
$ws.Range("M3").Copy()
$ws.Range("N3:O3").PasteSpecial(-4122)

# Row 4: year headers, same style as D4:M4 (s=31)
$ws.Range("M4").Copy()
$ws.Range("N4:O4").PasteSpecial(-4122)
$ws.Cells.Item(4, 14).Value = 2021
$ws.Cells.Item(4, 15).Value = 2022

# Row 5 ("Kyrgyz Republic" - bold row): copy format from M5, bump font size 9 -> 10
$ws.Range("M5").Copy()
$ws.Range("N5:O5").PasteSpecial(-4122)
$ws.Cells.Item(5, 14).Value = 40.007977647471066
$ws.Cells.Item(5, 15).Value = 42.620582506455563
$ws.Range("N5:O5").Font.Size = 10

# Rows 6-13 (oblast rows): copy format from M6:M13, bump font size 9 -> 10
$ws.Range("M6:M13").Copy()
$ws.Range("N6:N13").PasteSpecial(-4122)
$ws.Range("M6:M13").Copy()
$ws.Range("O6:O13").PasteSpecial(-4122)

$ws.Cells.Item(6, 14).Value = 5.7072514621689896
$ws.Cells.Item(6, 15).Value = 8.1443914479075037

$ws.Cells.Item(7, 14).Value = 8.9893229854028949
$ws.Cells.Item(7, 15).Value = 10.715961386284755

$ws.Cells.Item(8, 14).Value = 66.307512472824584
$ws.Cells.Item(8, 15).Value = 81.977461999426666

$ws.Cells.Item(9, 14).Value = 23.475213049310256
$ws.Cells.Item(9, 15).Value = 29.828871240443185

$ws.Cells.Item(10, 14).Value = 9.8045372040896162
$ws.Cells.Item(10, 15).Value = 9.7218425128664112

$ws.Cells.Item(11, 14).Value = 9.3737779268960448
$ws.Cells.Item(11, 15).Value = 8.6167819403064012

$ws.Cells.Item(12, 14).Value = 70.457032471318783
$ws.Cells.Item(12, 15).Value = 69.915337594090886

$ws.Cells.Item(13, 14).Value = 98.411252120183207
$ws.Cells.Item(13, 15).Value = 99.08571752721997

$ws.Range("N6:O13").Font.Size = 10

# Row 14 (total row with bottom border): copy format from M14, bump font size 9 -> 10
$ws.Range("M14").Copy()
$ws.Range("N14:O14").PasteSpecial(-4122)
$ws.Cells.Item(14, 14).Value = 63.900563564170795
$ws.Cells.Item(14, 15).Value = 64.805252627098838
$ws.Range("N14:O14").Font.Size = 10

# Move the active selection as in the source workbook
$ws.Range("P8").Select()
